$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23, shifting existing rows 23-27 down to 24-28.
$ws.Rows("23:23").Insert()

# Populate the newly inserted row 23 with a new Chirimoya price record
# (same boilerplate columns as the surrounding rows, new Fecha/Volumen).
$ws.Range("A23").Value = 5
$ws.Range("B23").Value = "Macroferia Regional de Talca"
$ws.Range("C23").Value = "Maule"
$ws.Range("D23").Value = "2021-09-29"
$ws.Range("E23").Value = 7
$ws.Range("F23").Value = "Fruta"
$ws.Range("G23").Value = 100107
$ws.Range("H23").Value = "Otros"
$ws.Range("I23").Value = 100107002
$ws.Range("J23").Value = "Chirimoya"
$ws.Range("K23").Value = "Cultivar IV Región"
$ws.Range("L23").Value = "Especial"
$ws.Range("M23").Value = 250
$ws.Range("N23").Value = 30000
$ws.Range("O23").Value = 30000
$ws.Range("P23").Value = 30000
$ws.Range("Q23").Value = "$/bandeja 10 kilos"
$ws.Range("R23").Value = "Provincia de Limarí"
$ws.Range("S23").Value = 3000
$ws.Range("T23").Value = 10

# Make sure the date cell keeps the same date/time number format used by
# the other "Fecha" cells in column D.
$ws.Range("D23").NumberFormat = $ws.Range("D24").NumberFormat
